$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.857.83"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.215.98"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.82"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.06%  "

$ws.Range("E9").Value = "  -1.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.29"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.00"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.547.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.21%  "

$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.224.79"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.783"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.827.77"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.06"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("E21").Value = "  -0.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.93%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "42.29"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.77"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.34"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.75%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.45"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0867"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +8.05%  "

$ws.Range("E34").Value = "  -1.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.122"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0356"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +7.93%  "

$ws.Range("E37").Value = "  -3.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.33"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.10"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.85"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +17.91%  "

$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "61.23"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.71%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.201"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.17%  "

$ws.Range("E44").Value = "  -2.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.491"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.54%  "

$ws.Range("E48").Value = "  -1.87%  "

$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.14"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.92%  "

$ws.Range("E51").Value = "  +19.34%  "
